# Auto-generated edit script: applies the 2024-10-18 YTD crime-data update
# across the Citywide Totals, By Neighborhood, and per-neighborhood sheets.
$wb = $excel.ActiveWorkbook

# --- Citywide Totals ---
$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("K2").Value = 6449
$ws.Range("K3").Value = 6649
$ws.Range("D4").Value = 1616
$ws.Range("F4").Value = 1582
$ws.Range("K4").Value = 1384
$ws.Range("K5").Value = 476
$ws.Range("K6").Value = 7325
$ws.Range("D7").Value = 22520
$ws.Range("F7").Value = 19404
$ws.Range("K7").Value = 22283

# --- Austin ---
$ws = $wb.Worksheets.Item("Austin")
$ws.Range("K2").Value = 403
$ws.Range("K3").Value = 445
$ws.Range("K4").Value = 83
$ws.Range("K6").Value = 490
$ws.Range("K7").Value = 1465

# --- South Chicago ---
$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("K2").Value = 164
$ws.Range("K3").Value = 171
$ws.Range("K4").Value = 24
$ws.Range("K7").Value = 479

# --- Garfield Park ---
$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("K2").Value = 250
$ws.Range("K7").Value = 975

# --- West Pullman ---
$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("K2").Value = 125
$ws.Range("K6").Value = 84
$ws.Range("K7").Value = 366

# --- Grand Crossing ---
$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("K2").Value = 215
$ws.Range("K6").Value = 226
$ws.Range("K7").Value = 756

# --- New City ---
$ws = $wb.Worksheets.Item("New City")
$ws.Range("K5").Value = 13
$ws.Range("K6").Value = 188
$ws.Range("K7").Value = 522

# --- Woodlawn ---
$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("K2").Value = 96
$ws.Range("K6").Value = 94
$ws.Range("K7").Value = 370

# --- By Neighborhood ---
$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("K2").Value = 195
$ws.Range("K7").Value = 666
$ws.Range("K8").Value = 1465
$ws.Range("K9").Value = 100
$ws.Range("K10").Value = 132
$ws.Range("K18").Value = 147
$ws.Range("K19").Value = 651
$ws.Range("K20").Value = 533
$ws.Range("K22").Value = 68
$ws.Range("K23").Value = 222
$ws.Range("K29").Value = 1202
$ws.Range("K31").Value = 247
$ws.Range("K33").Value = 975
$ws.Range("K37").Value = 756
$ws.Range("K42").Value = 823
$ws.Range("K43").Value = 181
$ws.Range("K44").Value = 185
$ws.Range("K46").Value = 45
$ws.Range("K48").Value = 279
$ws.Range("K51").Value = 283
$ws.Range("K54").Value = 438
$ws.Range("K55").Value = 242
$ws.Range("K56").Value = 24
$ws.Range("D63").Value = 269
$ws.Range("F63").Value = 161
$ws.Range("K63").Value = 58
$ws.Range("K65").Value = 522
$ws.Range("K67").Value = 872
$ws.Range("K72").Value = 115
$ws.Range("K74").Value = 24
$ws.Range("K76").Value = 304
$ws.Range("K77").Value = 153
$ws.Range("K79").Value = 559
$ws.Range("K82").Value = 24
$ws.Range("K83").Value = 479
$ws.Range("K84").Value = 179
$ws.Range("K85").Value = 1036
$ws.Range("K88").Value = 237
$ws.Range("K90").Value = 208
$ws.Range("K91").Value = 259
$ws.Range("K94").Value = 297
$ws.Range("K95").Value = 366
$ws.Range("K96").Value = 238
$ws.Range("K99").Value = 370
$ws.Range("D101").Value = 22520
$ws.Range("F101").Value = 19404
$ws.Range("K101").Value = 22283

# --- Gage Park ---
$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("K6").Value = 85
$ws.Range("K7").Value = 247

# --- North Lawndale ---
$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("K3").Value = 318
$ws.Range("K5").Value = 21
$ws.Range("K7").Value = 872

# --- South Deering ---
$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("K3").Value = 71
$ws.Range("K7").Value = 179

# --- Loop ---
$ws = $wb.Worksheets.Item("Loop")
$ws.Range("K2").Value = 71
$ws.Range("K6").Value = 236
$ws.Range("K7").Value = 438

# --- Englewood ---
$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("K2").Value = 341
$ws.Range("K6").Value = 346
$ws.Range("K7").Value = 1202

# --- Lake View ---
$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("K2").Value = 42
$ws.Range("K4").Value = 40
$ws.Range("K7").Value = 279

# --- Chatham ---
$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("K3").Value = 195
$ws.Range("K4").Value = 31
$ws.Range("K6").Value = 214
$ws.Range("K7").Value = 651

# --- Irving Park ---
$ws = $wb.Worksheets.Item("Irving Park")
$ws.Range("K6").Value = 74
$ws.Range("K7").Value = 185

# --- River North ---
$ws = $wb.Worksheets.Item("River North")
$ws.Range("K6").Value = 156
$ws.Range("K7").Value = 304

# --- Humboldt Park ---
$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("K2").Value = 224
$ws.Range("K6").Value = 305
$ws.Range("K7").Value = 823

# --- Avondale ---
$ws = $wb.Worksheets.Item("Avondale")
$ws.Range("K2").Value = 41
$ws.Range("K7").Value = 132

# --- Lower West Side ---
$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("K6").Value = 84
$ws.Range("K7").Value = 242

# --- Jefferson Park ---
$ws = $wb.Worksheets.Item("Jefferson Park")
$ws.Range("K5").Value = 1
$ws.Range("K7").Value = 45

# --- Douglas ---
$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("K2").Value = 64
$ws.Range("K7").Value = 222

# --- West Ridge ---
$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("K3").Value = 46
$ws.Range("K6").Value = 101
$ws.Range("K7").Value = 238

# --- Washington Park ---
$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("K3").Value = 122
$ws.Range("K7").Value = 259

# --- Roseland ---
$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("K2").Value = 187
$ws.Range("K7").Value = 559

# --- Chicago Lawn ---
$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("K2").Value = 184
$ws.Range("K6").Value = 145
$ws.Range("K7").Value = 533

# --- Calumet Heights ---
$ws = $wb.Worksheets.Item("Calumet Heights")
$ws.Range("K2").Value = 41
$ws.Range("K7").Value = 147

# --- Auburn Gresham ---
$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("K2").Value = 217
$ws.Range("K3").Value = 220
$ws.Range("K7").Value = 666

# --- West Loop ---
$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("K6").Value = 134
$ws.Range("K7").Value = 297

# --- Avalon Park ---
$ws = $wb.Worksheets.Item("Avalon Park")
$ws.Range("K2").Value = 32
$ws.Range("K7").Value = 100

# --- Albany Park ---
$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("K3").Value = 53
$ws.Range("K4").Value = 18
$ws.Range("K7").Value = 195

# --- United Center ---
$ws = $wb.Worksheets.Item("United Center")
$ws.Range("K3").Value = 72
$ws.Range("K7").Value = 237

# --- Washington Heights ---
$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("K6").Value = 51
$ws.Range("K7").Value = 208

# --- Little Italy, UIC ---
$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("K6").Value = 95
$ws.Range("K7").Value = 283

# --- Hyde Park ---
$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("K3").Value = 48
$ws.Range("K7").Value = 181

# --- South Shore ---
$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("K2").Value = 337
$ws.Range("K3").Value = 360
$ws.Range("K6").Value = 253
$ws.Range("K7").Value = 1036

# --- Clearing ---
$ws = $wb.Worksheets.Item("Clearing")
$ws.Range("K2").Value = 32
$ws.Range("K7").Value = 68

# --- Old Town ---
$ws = $wb.Worksheets.Item("Old Town")
$ws.Range("K6").Value = 55
$ws.Range("K7").Value = 115

# --- Sheffield & DePaul ---
$ws = $wb.Worksheets.Item("Sheffield & DePaul")
$ws.Range("K4").Value = 4
$ws.Range("K6").Value = 24

# --- Riverdale ---
$ws = $wb.Worksheets.Item("Riverdale")
$ws.Range("K4").Value = 12
$ws.Range("K7").Value = 153

# --- Magnificent Mile ---
$ws = $wb.Worksheets.Item("Magnificent Mile")
$ws.Range("K4").Value = 2
$ws.Range("K7").Value = 24

# --- Printers Row ---
$ws = $wb.Worksheets.Item("Printers Row")
$ws.Range("K6").Value = 15
$ws.Range("K7").Value = 24
